$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preventing Excel from
# auto-converting numeric-looking strings (e.g. "1.003") into numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Updated Price (D) / Volume(1h) (E) figures for the crypto list refresh.
$updates = @(
    @{ Row = 2; D = "27.454.08"; E = "  -1.94%  " },
    @{ Row = 3; D = "1.831.29"; E = "  -2.60%  " },
    @{ Row = 4; D = "1.003"; E = "  -0.71%  " },
    @{ Row = 5; D = "330.88"; E = "  -1.65%  " },
    @{ Row = 6; D = "1.003"; E = "  -0.75%  " },
    @{ Row = 7; D = "0.4593"; E = "  -3.84%  " },
    @{ Row = 8; D = "0.3817"; E = "  -3.34%  " },
    @{ Row = 9; D = "46.70"; E = "  -0.87%  " },
    @{ Row = 10; D = "0.07906"; E = "  -1.44%  " },
    @{ Row = 11; D = "0.9691"; E = "  -4.91%  " },
    @{ Row = 12; D = "21.00"; E = "  -4.12%  " },
    @{ Row = 13; D = "1.832.81"; E = "  -2.96%  " },
    @{ Row = 14; D = "5.885"; E = "  -2.86%  " },
    @{ Row = 15; D = "7.036"; E = "  -2.35%  " },
    @{ Row = 16; D = $null; E = "  -0.86%  " },
    @{ Row = 17; D = "87.91"; E = "  -0.94%  " },
    @{ Row = 18; D = "0.06620"; E = "  -1.87%  " },
    @{ Row = 19; D = "0.00001028"; E = "  -2.22%  " },
    @{ Row = 20; D = "17.00"; E = "  -0.50%  " },
    @{ Row = 21; D = "1.003"; E = "  -0.64%  " },
    @{ Row = 22; D = "27.453.02"; E = "  -1.88%  " },
    @{ Row = 23; D = "5.345"; E = "  -3.04%  " },
    @{ Row = 24; D = $null; E = "  -1.94%  " },
    @{ Row = 25; D = "2.309"; E = "  -1.47%  " },
    @{ Row = 26; D = "2.050.91"; E = "  -2.73%  " },
    @{ Row = 27; D = $null; E = "  -0.72%  " },
    @{ Row = 28; D = "19.35"; E = "  -2.77%  " },
    @{ Row = 29; D = $null; E = "  -2.36%  " },
    @{ Row = 30; D = "5.297"; E = $null },
    @{ Row = 31; D = "118.93"; E = "  -2.17%  " },
    @{ Row = 32; D = "0.9531"; E = "  -2.72%  " },
    @{ Row = 33; D = "0.09294"; E = "  -2.76%  " },
    @{ Row = 34; D = "3.583"; E = "  -1.46%  " },
    @{ Row = 35; D = "5.240"; E = "  -1.98%  " },
    @{ Row = 36; D = "1.311"; E = "  -3.37%  " },
    @{ Row = 37; D = "0.05924"; E = "  -2.74%  " },
    @{ Row = 38; D = "0.02192"; E = "  -2.41%  " },
    @{ Row = 39; D = "1.162"; E = "  -3.88%  " },
    @{ Row = 40; D = "8.041"; E = "  -2.08%  " },
    @{ Row = 41; D = "0.5781"; E = "  -3.48%  " },
    @{ Row = 42; D = "0.1836"; E = "  -3.33%  " },
    @{ Row = 43; D = "10.03"; E = "  -3.23%  " },
    @{ Row = 44; D = "1.259"; E = "  +0.05%  " },
    @{ Row = 45; D = "0.5477"; E = "  -3.66%  " },
    @{ Row = 46; D = "11.92"; E = "  -2.36%  " },
    @{ Row = 47; D = "1.862"; E = "  -3.81%  " },
    @{ Row = 48; D = $null; E = "  -2.37%  " },
    @{ Row = 49; D = "110.28"; E = "  -2.21%  " },
    @{ Row = 50; D = "1.038"; E = "  -3.17%  " },
    @{ Row = 51; D = "1.002"; E = "  -0.91%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $isNumericLooking = $u.D -match '^[0-9]+(\.[0-9]+)?$'
        if ($isNumericLooking) {
            Set-TextValue $ws.Range("D" + $u.Row) $u.D
        } else {
            $ws.Range("D" + $u.Row).Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
